# Update currency quotation values and "last updated" timestamps
# in the "moedas_atualizadas" worksheet, per the latest maestro run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => [Column E new value (optional), Column F new value]
$updates = @(
    @{ Row = 2;  E = "5,57"; F = "11 de jun., 12:13 UTC ·" },
    @{ Row = 3;  E = "6,37"; F = "11 de jun., 12:13 UTC ·" },
    @{ Row = 4;  E = $null; F = "11 de jun., 12:13 UTC ·" },
    @{ Row = 5;  E = "7,52"; F = "11 de jun., 12:13 UTC ·" },
    @{ Row = 6;  E = $null; F = "11 de jun., 12:13 UTC ·" },
    @{ Row = 7;  E = "6,77"; F = "11 de jun., 12:13 UTC ·" },
    @{ Row = 8;  E = "3,62"; F = "11 de jun., 12:13 UTC ·" },
    @{ Row = 9;  E = $null; F = "11 de jun., 12:13 UTC ·" },
    @{ Row = 10; E = $null; F = "11 de jun., 12:12 UTC ·" },
    @{ Row = 11; E = $null; F = "11 de jun., 12:13 UTC ·" },
    @{ Row = 12; E = "0,78"; F = "11 de jun., 12:14 UTC ·" },
    @{ Row = 13; E = $null; F = "11 de jun., 12:12 UTC ·" },
    @{ Row = 14; E = $null; F = "11 de jun., 11:19 UTC ·" },
    @{ Row = 15; E = $null; F = "11 de jun., 12:09 UTC ·" },
    @{ Row = 16; E = $null; F = "11 de jun., 12:13 UTC ·" },
    @{ Row = 17; E = "1,07"; F = "11 de jun., 12:13 UTC ·" },
    @{ Row = 18; E = $null; F = "11 de jun., 12:13 UTC ·" },
    @{ Row = 19; E = "4,33"; F = "11 de jun., 12:13 UTC ·" },
    @{ Row = 20; E = $null; F = "11 de jun., 12:13 UTC ·" },
    @{ Row = 21; E = "4,16"; F = "11 de jun., 12:13 UTC ·" },
    @{ Row = 22; E = $null; F = "11 de jun., 12:13 UTC ·" },
    @{ Row = 23; E = "3,50"; F = "11 de jun., 12:13 UTC ·" },
    @{ Row = 24; E = "0,78"; F = "11 de jun., 12:14 UTC ·" }
)

foreach ($u in $updates) {
    if ($null -ne $u.E) {
        $ws.Range("E$($u.Row)").Value = $u.E
    }
    $ws.Range("F$($u.Row)").Value = $u.F
}
